$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "regular"
$ws.Range("B23").Value = "full random"
$ws.Range("C23").Value = "sentences"
$ws.Range("D23").Value = 5000
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 14000
$ws.Range("I23").Value = "auto"
$ws.Range("J23").Value = "yes"
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = "bad"

$ws.Range("A24").Value = "zeroshot huang combined with own (w/o marketing)"
$ws.Range("D24").Value = 4000
$ws.Range("E24").Value = 200
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 0.05
$ws.Range("H24").Value = 4000
$ws.Range("I24").Value = "null"
$ws.Range("J24").Value = "yes"
$ws.Range("K24").Value = 14
$ws.Range("L24").Value = "representations very good; but only few outliers, outliers might be assigned to other topics. This needs check in the full workflow. "

$ws.Range("A24:L24").Interior.Color = 65535

$ws.Range("G26").Select()
